$d = $word.ActiveDocument

# Locate the "Metaclass, Class, Instance, Context, Occurrence, Role Resource
# Metadata Maps Monad..." bullet -- the new content is inserted right after it
# (and before the blank paragraph that already precedes "Java pattern
# matching...").
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Metaclass, Class, Instance, Context, Occurrence, Role Resource Metadata Maps Monad with contextual CSPOs Statements*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph"
}

$anchor = $d.Paragraphs($anchorIndex)

# New bullet (numId=3) list items to add, in document order.
$items = @(
    "Resources: John, Peter, Mary, loves, friendOf, loverHasFriend.",
    "Transform / Mapping: John :loves Mary;",
    "Transform / Mapping: Peter :friendOf John;",
    "Transform / Mapping (Expanded Knowledge): Mary :loverHasFriend Peter;",
    "Browse Resources: Functional Activation: Transforms / Mappings Contexts Knowledge Expansion. Idem for Kinds and Schema Statements Aggregations / Alignments."
)

# $anchor.Range.InsertParagraphAfter() always inserts a new paragraph
# immediately after $anchor (inheriting $anchor's pPr -- i.e. the numId=3
# bullet formatting), so walking the list back-to-front yields the paragraphs
# in the correct final reading order.
for ($j = $items.Length - 1; $j -ge 0; $j--) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($anchorIndex + 1)
    $newPara.Range.Text = $items[$j]
    $newPara.Range.Font.Underline = 0
}

# Finally, insert the new leading blank paragraph right after $anchor (so it
# ends up ahead of the five bullets above) and strip the inherited bullet
# formatting so it matches the plain blank paragraph that already sits just
# below the new block.
$anchor.Range.InsertParagraphAfter()
$blank = $d.Paragraphs($anchorIndex + 1)
$blank.Range.ListFormat.RemoveNumbers()
$blank.LeftIndent = 0
$blank.FirstLineIndent = 0

Write-Host "Paragraphs after edit: $($d.Paragraphs.Count)"
